$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row (row 1: "Product","Reactant","Sentences") and shift
# all data rows up by one.
$ws.Rows.Item(1).Delete()

# Reset the view to match the updated sheet (top-left cell A1, selection D6).
$ws.Range("D6").Select()
